$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.074.50"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "1.800.04"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "308.08"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4214"
$ws.Range("E7").Value = "  -2.36%  "
$ws.Range("D8").Value = "0.3602"
$ws.Range("D9").Value = "0.07219"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "0.8472"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").Value = "20.27"
$ws.Range("E11").Value = "  -3.63%  "
$ws.Range("D12").Value = "1.887.55"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("D13").Value = "5.304"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "6.379"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "0.06761"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "80.39"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "0.000008717"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "15.07"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").Value = "27.015.69"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("D22").Value = "5.081"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "2.011.41"
$ws.Range("E24").Value = "  -7.33%  "
$ws.Range("D25").Value = "1.927"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").Value = "153.15"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "18.20"
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("D28").Value = "5.038"
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("D29").Value = "113.54"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "1.658"
$ws.Range("E30").Value = "  -11.96%  "
$ws.Range("D31").Value = "0.09009"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "0.7292"
$ws.Range("E32").Value = "  -7.21%  "
$ws.Range("D33").Value = "2.862"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "4.340"
$ws.Range("D35").Value = "1.097"
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "1.080"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "0.05152"
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("D39").Value = "0.01907"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").Value = "0.4982"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").Value = "2.619"
$ws.Range("E42").Value = "  -7.82%  "
$ws.Range("D43").Value = "8.090"
$ws.Range("E43").Value = "  -6.64%  "
$ws.Range("D44").Value = "5.983"
$ws.Range("E44").Value = "  -11.92%  "
$ws.Range("D45").Value = "105.29"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "10.19"
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "0.06301"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").Value = "0.4550"
$ws.Range("E49").Value = "  -5.17%  "
$ws.Range("D50").Value = "1.607"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Value = "1.726"
$ws.Range("E51").Value = "  -6.52%  "
